# Replace numeric AS/NZS1170.2 "version" literals (2011 / 2021) in column A
# of several lookup sheets with descriptive enum-style string labels, and
# restore the selection / active-sheet state recorded in the edited file.

$wb = $excel.ActiveWorkbook

$v2011 = "AS/NZS1170.2-2011"
$v2021 = "AS/NZS1170.2-2021"

# --- shielding_multiplier ---------------------------------------------
$ws = $wb.Worksheets.Item("shielding_multiplier")
$ws.Range("A2:A5").Value = $v2011
$ws.Range("A6:A9").Value = $v2021

# --- region_windspeed_parameters ---------------------------------------
$ws = $wb.Worksheets.Item("region_windspeed_parameters")
$ws.Range("A2:A12").Value = $v2011
$ws.Range("A13:A26").Value = $v2021

# --- region_direction_parameters ---------------------------------------
$ws = $wb.Worksheets.Item("region_direction_parameters")
$ws.Range("A2:A111").Value = $v2011
$ws.Range("A112:A251").Value = $v2021

# --- terrain_height_multipliers -----------------------------------------
$ws = $wb.Worksheets.Item("terrain_height_multipliers")
$ws.Range("A2:A53").Value = $v2011
$ws.Range("A54:A118").Value = $v2021

# --- cpi_t5a --------------------------------------------------------------
$ws = $wb.Worksheets.Item("cpi_t5a")
$ws.Range("A2:A9").Value = $v2021

# --- cpi_t5b --------------------------------------------------------------
$ws = $wb.Worksheets.Item("cpi_t5b")
$ws.Range("A2:A21").Value = $v2021

# --- k_a -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("k_a")
$ws.Range("A2:A13").Value = $v2021

# --- cpe_t5_2c ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("cpe_t5_2c")
$ws.Range("A2:A7").Value = $v2021
$ws.Range("E4").Select()

# --- wind_direction_definitions: just a cursor move ------------------------
$ws = $wb.Worksheets.Item("wind_direction_definitions")
$ws.Range("D8").Select()

# --- region_direction_parameters: cursor move -------------------------------
$ws = $wb.Worksheets.Item("region_direction_parameters")
$ws.Range("G6").Select()

# --- app_c_fig_c2: cursor move (was the active tab before the edit) --------
$ws = $wb.Worksheets.Item("app_c_fig_c2")
$ws.Range("H8").Select()

# --- shielding_multiplier becomes the active tab, cursor on D9 -------------
$ws = $wb.Worksheets.Item("shielding_multiplier")
$ws.Range("D9").Select()
$ws.Activate()
